$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023-05 Compr. Table Linux")

# Row 17 updates (ROW STORE COMPRESS BASIC)
$ws.Range("C17").Value = 60
$ws.Range("J17").Value = 13.744
$ws.Range("K17").Value = 18.977
$ws.Range("L17").Value = 0.0029
$ws.Range("M17").Value = 0.609
$ws.Range("N17").Value = 20812
$ws.Range("O17").Value = 0.344
$ws.Range("P17").Value = 0.633
$ws.Range("Q17").Value = 21137

# Row 18 updates (ROW STORE COMPRESS ADVANCED) - previously empty
$ws.Range("B18").Value = 168
$ws.Range("C18").Value = 60
$ws.Range("D18").Value = 0.297
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 16.956
$ws.Range("G18").Value = 0.131
$ws.Range("H18").Value = 1.422
$ws.Range("I18").Value = 0.131
$ws.Range("J18").Value = 14.182
$ws.Range("K18").Value = 19.933
$ws.Range("L18").Value = 0.0029
$ws.Range("M18").Value = 0.626
$ws.Range("N18").Value = 20812
$ws.Range("O18").Value = 0.382
$ws.Range("P18").Value = 0.717
$ws.Range("Q18").Value = 21115

# Update the active cell selection on this sheet
$ws.Range("P21").Select()
